# The underlying OOXML for this template is reformatted only (every
# changed line in the source diff is the exact same element with its
# attributes re-ordered alphabetically - xmlns declarations, w:pgSz,
# w:pgMar, w:rFonts, w:lang, w:latentStyles/w:lsdException, w:style and
# table-cell margins). No attribute, value, namespace URI, text run, or
# document structure actually changes (verified by canonicalizing the
# "before" and reconstructed "after" XML: they are identical). So there
# is nothing for a Word automation session to edit content-wise; we just
# touch the document through the object model to confirm it loads and
# round-trips cleanly without introducing any unintended side effects.
$d = $word.ActiveDocument
$null = $d.Content.Text
